$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 5.21898235657733
$ws.Range("D2").Value = 0.3375769260293227
$ws.Range("E2").Value = 21.511472964381
$ws.Range("F2").Value = 4.517967691790099
$ws.Range("G2").Value = 5.919997021364561

$ws.Range("C3").Value = 7.464574684069104
$ws.Range("D3").Value = 0.3375769260293215
$ws.Range("E3").Value = 21.51147296438204
$ws.Range("F3").Value = 6.763560019281877
$ws.Range("G3").Value = 8.165589348856331

$ws.Range("C4").Value = 5.27558953566662
$ws.Range("D4").Value = 0.6000400152698083
$ws.Range("E4").Value = 16.44637937837114
$ws.Range("F4").Value = 4.006361686860882
$ws.Range("G4").Value = 6.544817384472357

$ws.Range("C5").Value = 7.790703678333328
$ws.Range("D5").Value = 0.6000400152698085
$ws.Range("E5").Value = 16.44637937837151
$ws.Range("F5").Value = 6.521475829527592
$ws.Range("G5").Value = 9.059931527139064

$ws.Range("C6").Value = 4.175905373052956
$ws.Range("D6").Value = 0.7532527608297563
$ws.Range("E6").Value = 23.40073583718799
$ws.Range("F6").Value = 2.619158544809541
$ws.Range("G6").Value = 5.732652201296371

$ws.Range("C7").Value = 8.125867857332864
$ws.Range("D7").Value = 0.7532527608297541
$ws.Range("E7").Value = 23.40073583718769
$ws.Range("F7").Value = 6.569121029089453
$ws.Range("G7").Value = 9.682614685576276

$ws.Range("C8").Value = 4.890159088432354
$ws.Range("D8").Value = 0.3726900639192939
$ws.Range("E8").Value = 59.09356675379004
$ws.Range("F8").Value = 4.144432698497401
$ws.Range("G8").Value = 5.635885478367308

$ws.Range("C9").Value = 7.793715406578436
$ws.Range("D9").Value = 0.3726900639192933
$ws.Range("E9").Value = 59.09356675379071
$ws.Range("F9").Value = 7.047989016643485
$ws.Range("G9").Value = 8.539441796513387
